$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed crypto price/volume snapshot values.
# Column D ("Price") values such as "75.944.06" or "1.20" look numeric/
# date-like, so force a text number format on those specific cells before
# assigning them -- otherwise Excel auto-converts them to numbers and
# trailing zeros / grouping dots would be lost.

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '75.944.06'

$ws.Cells.Item(2, 5).Value = '  +1.62%  '

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '2.908.15'

$ws.Cells.Item(3, 5).Value = '  +2.51%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '197.90'

$ws.Cells.Item(5, 5).Value = '  +4.82%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '596.34'

$ws.Cells.Item(6, 5).Value = '  -0.63%  '

$ws.Cells.Item(7, 5).Value = '  -0.04%  '

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = '@'
$c.Value = '0.550'

$ws.Cells.Item(8, 5).Value = '  -1.11%  '

$ws.Cells.Item(9, 5).Value = '  +0.67%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '2.907.79'

$ws.Cells.Item(10, 5).Value = '  +2.44%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.420'

$ws.Cells.Item(11, 5).Value = '  +13.16%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '0.160'

$ws.Cells.Item(12, 5).Value = '  -1.14%  '

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '4.88'

$ws.Cells.Item(13, 5).Value = '  -0.54%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '3.440.97'

$ws.Cells.Item(14, 5).Value = '  +2.29%  '

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '75.776.97'

$ws.Cells.Item(15, 5).Value = '  +1.21%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '0.0000189'

$ws.Cells.Item(16, 5).Value = '  +0.42%  '

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '27.28'

$ws.Cells.Item(17, 5).Value = '  -0.09%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '2.903.86'

$ws.Cells.Item(18, 5).Value = '  +1.75%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = '8.82'

$ws.Cells.Item(19, 5).Value = '  -3.54%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '12.74'

$ws.Cells.Item(20, 5).Value = '  +2.48%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '376.76'

$ws.Cells.Item(21, 5).Value = '  +0.37%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '2.30'

$ws.Cells.Item(22, 5).Value = '  +0.99%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '4.17'

$ws.Cells.Item(23, 5).Value = '  +1.21%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = '71.24'

$ws.Cells.Item(25, 5).Value = '  +0.05%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '3.059.77'

$ws.Cells.Item(26, 5).Value = '  +2.31%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '4.18'

$ws.Cells.Item(27, 5).Value = '  -0.88%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '9.61'

$ws.Cells.Item(28, 5).Value = '  -0.07%  '

$ws.Cells.Item(29, 5).Value = '  +4.76%  '

$ws.Cells.Item(30, 5).Value = '  -0.10%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '1.40'

$ws.Cells.Item(31, 5).Value = '  -0.10%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = '@'
$c.Value = '501.05'

$ws.Cells.Item(32, 5).Value = '  -4.88%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '7.71'

$ws.Cells.Item(33, 5).Value = '  -2.32%  '

$ws.Cells.Item(34, 5).Value = '  -0.70%  '

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '1.00'

$ws.Cells.Item(35, 5).Value = '  +0.02%  '

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '164.12'

$ws.Cells.Item(36, 5).Value = '  +1.21%  '

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '20.06'

$ws.Cells.Item(37, 5).Value = '  +0.01%  '

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '19.69'

$ws.Cells.Item(39, 5).Value = '  -5.98%  '

$ws.Cells.Item(40, 5).Value = '  -0.07%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '180.20'

$ws.Cells.Item(41, 5).Value = '  -0.87%  '

$ws.Cells.Item(42, 5).Value = '  +0.70%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '4.98'

$ws.Cells.Item(43, 5).Value = '  -2.08%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '1.65'

$ws.Cells.Item(44, 5).Value = '  -2.52%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.0907'

$ws.Cells.Item(45, 5).Value = '  +6.85%  '

$ws.Cells.Item(46, 2).Value = 'OKB'

$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '40.11'

$ws.Cells.Item(46, 5).Value = '  +0.95%  '

$ws.Cells.Item(47, 2).Value = 'ImmutableX'

$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '1.20'

$ws.Cells.Item(47, 5).Value = '  -3.85%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '2.31'

$ws.Cells.Item(48, 5).Value = '  -2.57%  '

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '0.575'

$ws.Cells.Item(49, 5).Value = '  +0.72%  '

$ws.Cells.Item(50, 2).Value = 'Mantle'

$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '0.660'

$ws.Cells.Item(50, 5).Value = '  +6.80%  '

$ws.Cells.Item(51, 2).Value = 'Filecoin'

$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '3.71'

$ws.Cells.Item(51, 5).Value = '  -1.03%  '
